$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.787.63"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "'3.439.21"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'573.25"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").Value = "'158.19"
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.605"
$ws.Range("E8").Value = "  +3.24%  "
$ws.Range("D9").Value = "'3.439.09"
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").Value = "'7.26"
$ws.Range("E10").Value = "  -1.68%  "
$ws.Range("D11").Value = "'0.123"
$ws.Range("E11").Value = "  -1.89%  "
$ws.Range("D12").Value = "'0.447"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "'4.036.69"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("D14").Value = "'0.134"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("E15").Value = "  -2.57%  "
$ws.Range("D16").Value = "'28.11"
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("D17").Value = "'64.824.46"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "'3.458.49"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").Value = "'6.39"
$ws.Range("E19").Value = "  -1.20%  "
$ws.Range("D20").Value = "'14.16"
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("D21").Value = "'376.56"
$ws.Range("E21").Value = "  -4.18%  "
$ws.Range("D22").Value = "'8.08"
$ws.Range("E22").Value = "  -2.84%  "
$ws.Range("D23").Value = "'0.553"
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("D24").Value = "'72.50"
$ws.Range("E24").Value = "  -1.77%  "
$ws.Range("D25").Value = "'0.997"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("E26").Value = "  -2.46%  "
$ws.Range("E27").Value = "  +4.52%  "
$ws.Range("D28").Value = "'0.177"
$ws.Range("E28").Value = "  -1.48%  "
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("D30").Value = "'1.49"
$ws.Range("E30").Value = "  +2.51%  "
$ws.Range("D31").Value = "'6.10"
$ws.Range("E31").Value = "  -1.94%  "
$ws.Range("D32").Value = "'2.04"
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("D33").Value = "'23.40"
$ws.Range("E33").Value = "  -1.54%  "
$ws.Range("D34").Value = "'7.24"
$ws.Range("E34").Value = "  +2.55%  "
$ws.Range("D35").Value = "'1.59"
$ws.Range("E35").Value = "  +6.74%  "
$ws.Range("D36").Value = "'159.39"
$ws.Range("E36").Value = "  -1.03%  "
$ws.Range("D37").Value = "'1.91"
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("D38").Value = "'0.0771"
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("D39").Value = "'27.41"
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").Value = "'4.72"
$ws.Range("E40").Value = "  +6.19%  "
$ws.Range("D41").Value = "'6.81"
$ws.Range("E41").Value = "  +2.87%  "
$ws.Range("D42").Value = "'2.867.74"
$ws.Range("E42").Value = "  -2.38%  "
$ws.Range("D43").Value = "'42.98"
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("D44").Value = "'0.0317"
$ws.Range("E44").Value = "  -1.02%  "
$ws.Range("D45").Value = "'0.775"
$ws.Range("E45").Value = "  +0.32%  "
$ws.Range("D46").Value = "'25.82"
$ws.Range("E46").Value = "  +8.62%  "
$ws.Range("D47").Value = "'320.76"
$ws.Range("E47").Value = "  +8.40%  "
$ws.Range("D48").Value = "'1.08"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("D49").Value = "'0.109"
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("D50").Value = "'0.867"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").Value = "'6.60"
$ws.Range("E51").Value = "  +0.63%  "
